$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Strip the old bold header formatting and the per-row border box ---
$ws.Range("A1:D1").Font.Bold = $false
$ws.Range("A1:D1").Borders.LineStyle = 0
$ws.Range("B2:D3").Borders.LineStyle = 0

# --- New header row (row 1) ---
$ws.Range("A1").Value = "name"
$ws.Range("B1").Value = "type"
$ws.Range("C1").Value = "robots"
$ws.Range("D1").Value = "propertiesTestSingleTask.createTask"

# --- Row 2 ---
$ws.Range("A2").Value = "SingleTaskScenario#TestForCreatingTodoList"
$ws.Range("B2").Value = "test"
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = "os:win;chrome:80"

# --- Row 3 ---
$ws.Range("A3").Value = "SingleTaskScenario#AnotherTest"
$ws.Range("B3").Value = "test"
$ws.Range("C3").ClearContents()
$ws.Range("D3").Value = "os:win;chrome:80"

# --- Column widths: column A widens to fit the longer names, column B reverts ---
$ws.Columns("B").ColumnWidth = 8.43
$ws.Columns("A").ColumnWidth = 40.5

# --- Selection moves to B13 ---
$ws.Range("B13").Select()
